$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "303.14"
    "E2" = "4.73%"
    "E3" = "9.68%"
    "D4" = "5.259"
    "E4" = "-0.12%"
    "D5" = "0.07521"
    "E5" = "6.28%"
    "D6" = "7.888"
    "E6" = "5.68%"
    "D7" = "3.817"
    "E7" = "6.97%"
    "D8" = "1.514"
    "E8" = "8.56%"
    "D9" = "0.9212"
    "E9" = "1.45%"
    "D10" = "0.1698"
    "E10" = "4.91%"
    "D11" = "0.07970"
    "E11" = "4.31%"
    "D12" = "0.08025"
    "E12" = "3.17%"
    "E13" = "3.87%"
    "D14" = "0.09916"
    "E14" = "9.75%"
    "D15" = "0.001495"
    "E15" = "-4.95%"
    "D16" = "0.04610"
    "E16" = "1.80%"
    "D17" = "0.006440"
    "E17" = "3.94%"
    "D18" = "3.460"
    "E18" = "-0.74%"
    "D19" = "2.226"
    "E19" = "-0.29%"
    "D20" = "0.3302"
    "E20" = "2.18%"
    "D21" = "0.1343"
    "E21" = "-0.50%"
    "D22" = "4.499"
    "E22" = "12.25%"
    "E23" = "1.24%"
    "D24" = "0.001218"
    "E24" = "0.70%"
    "D25" = "0.004447"
    "E25" = "6.22%"
    "D26" = "0.0001399"
    "E26" = "19.75%"
    "D27" = "0.0001774"
    "E27" = "5.99%"
    "D39" = "0.01719"
    "E39" = "2,542.94%"
    "D40" = "0.04489"
    "E40" = "2.33%"
    "D41" = "0.006958"
    "D42" = "0.1348"
    "E42" = "6.99%"
    "D43" = "0.002128"
    "E43" = "2.99%"
    "D44" = "0.01288"
    "E44" = "9.25%"
    "D45" = "0.00006164"
    "E45" = "5.59%"
    "D46" = "0.7094"
    "E46" = "-63.23%"
    "D47" = "0.01498"
    "E47" = "15.43%"
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.ClearFormats()
}
